# "In progress CLS+Thủ Thuật" -- bump the receive-record ids / insurance
# card / id-card numbers on rows 2 and 3 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("A2").Value = 487
$ws.Range("E2").Value = 46200021274
$ws.Range("X2").Value = "DN4127389127785"

# Row 3
$ws.Range("A3").Value = 488
$ws.Range("E3").Value = 46200021275

# Mirror the author's cursor/selection landing on the InsCardNo cell.
$ws.Range("X2").Select()
